# Adapt tests to control version
# Add a "version" column (header + value) to the settings sheet, and make
# the settings sheet the active tab (selection moves to C3).

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("settings")

# New header "version" in C1 and value 1 in C2
$settings.Cells.Item(1, 3).Value = "version"
$settings.Cells.Item(2, 3).Value = 1

# Make "settings" the active sheet/tab, and move the selection to C3
$settings.Activate()
$settings.Range("C3").Select()
